$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 5.019802666666667
$ws.Range("H2").Value = 15.059408
$ws.Range("I2").Value = 0.1084294745534388
$ws.Range("J2").Value = 0.1084294745534388
$ws.Range("M2").Value = 0.6692693333333334
$ws.Range("N2").Value = 2.007808
$ws.Range("O2").Value = 0.004126561180566838
$ws.Range("P2").Value = 0.004126561180566839
$ws.Range("Q2").Value = 3.359599984184889
$ws.Range("R2").Value = 30.23639985766401
$ws.Range("S2").Value = 0.0004474408605214803
$ws.Range("T2").Value = 0.0004474408605214804

$ws.Range("G3").Value = 5.019802666666667
$ws.Range("H3").Value = 15.059408
$ws.Range("I3").Value = 0.1084294745534388
$ws.Range("J3").Value = 0.1084294745534388
$ws.Range("O3").Value = 0.9916964991825307
$ws.Range("P3").Value = 0.9916964991825309
$ws.Range("Q3").Value = 807.3801398267857
$ws.Range("R3").Value = 7266.421258441072
$ws.Range("S3").Value = 0.1075291303228465
$ws.Range("T3").Value = 0.1075291303228466

$ws.Range("G4").Value = 5.019802666666667
$ws.Range("H4").Value = 15.059408
$ws.Range("I4").Value = 0.1084294745534388
$ws.Range("J4").Value = 0.1084294745534388
$ws.Range("M4").Value = 0.5637343333333333
$ws.Range("N4").Value = 1.691203
$ws.Range("O4").Value = 0.003475856580040611
$ws.Range("P4").Value = 0.003475856580040611
$ws.Range("Q4").Value = 2.829835109758223
$ws.Range("R4").Value = 25.468515987824
$ws.Range("S4").Value = 0.0003768853025969161
$ws.Range("T4").Value = 0.0003768853025969162

$ws.Range("G5").Value = 5.019802666666667
$ws.Range("H5").Value = 15.059408
$ws.Range("I5").Value = 0.1084294745534388
$ws.Range("J5").Value = 0.1084294745534388
$ws.Range("M5").Value = 0.1137056666666667
$ws.Range("N5").Value = 0.341117
$ws.Range("O5").Value = 0.0007010830568617209
$ws.Range("P5").Value = 0.0007010830568617211
$ws.Range("Q5").Value = 0.5707800087484445
$ws.Range("R5").Value = 5.137020078736001
$ws.Range("S5").Value = 0.00007601806747383503
$ws.Range("T5").Value = 0.00007601806747383506

$ws.Range("G6").Value = 4.734454666666666
$ws.Range("I6").Value = 0.1022658590172488
$ws.Range("J6").Value = 0.1022658590172488
$ws.Range("M6").Value = 0.6692693333333334
$ws.Range("N6").Value = 2.007808
$ws.Range("O6").Value = 0.004126561180566838
$ws.Range("P6").Value = 0.004126561180566839
$ws.Range("Q6").Value = 3.168625318456889
$ws.Range("R6").Value = 28.517627866112
$ws.Range("S6").Value = 0.0004220063239178999
$ws.Range("T6").Value = 0.0004220063239179

$ws.Range("G7").Value = 4.734454666666666
$ws.Range("I7").Value = 0.1022658590172488
$ws.Range("J7").Value = 0.1022658590172488
$ws.Range("O7").Value = 0.9916964991825307
$ws.Range("P7").Value = 0.9916964991825309
$ws.Range("Q7").Value = 761.4850472429416
$ws.Range("S7").Value = 0.1014166943732998
$ws.Range("T7").Value = 0.1014166943732999

$ws.Range("G8").Value = 4.734454666666666
$ws.Range("I8").Value = 0.1022658590172488
$ws.Range("J8").Value = 0.1022658590172488
$ws.Range("M8").Value = 0.5637343333333333
$ws.Range("N8").Value = 1.691203
$ws.Range("O8").Value = 0.003475856580040611
$ws.Range("P8").Value = 0.003475856580040611
$ws.Range("Q8").Value = 2.668974645210222
$ws.Range("R8").Value = 24.020771806892
$ws.Range("S8").Value = 0.0003554614589786095
$ws.Range("T8").Value = 0.0003554614589786096

$ws.Range("G9").Value = 4.734454666666666
$ws.Range("I9").Value = 0.1022658590172488
$ws.Range("J9").Value = 0.1022658590172488
$ws.Range("M9").Value = 0.1137056666666667
$ws.Range("N9").Value = 0.341117
$ws.Range("O9").Value = 0.0007010830568617209
$ws.Range("P9").Value = 0.0007010830568617211
$ws.Range("Q9").Value = 0.5383343241764443
$ws.Range("R9").Value = 4.845008917587999
$ws.Range("S9").Value = 0.00007169686105240254
$ws.Range("T9").Value = 0.00007169686105240257

$ws.Range("G10").Value = 28.73869333333333
$ws.Range("H10").Value = 86.21608000000001
$ws.Range("I10").Value = 0.6207657201702246
$ws.Range("J10").Value = 0.6207657201702246
$ws.Range("M10").Value = 0.6692693333333334
$ws.Range("N10").Value = 2.007808
$ws.Range("O10").Value = 0.004126561180566838
$ws.Range("P10").Value = 0.004126561180566839
$ws.Range("Q10").Value = 19.23392612807111
$ws.Range("R10").Value = 173.10533515264
$ws.Range("S10").Value = 0.002561627723081065
$ws.Range("T10").Value = 0.002561627723081066

$ws.Range("G11").Value = 28.73869333333333
$ws.Range("H11").Value = 86.21608000000001
$ws.Range("I11").Value = 0.6207657201702246
$ws.Range("J11").Value = 0.6207657201702246
$ws.Range("O11").Value = 0.9916964991825307
$ws.Range("P11").Value = 0.9916964991825309
$ws.Range("Q11").Value = 4622.303262234302
$ws.Range("R11").Value = 41600.72936010872
$ws.Range("S11").Value = 0.6156111915053342
$ws.Range("T11").Value = 0.6156111915053343

$ws.Range("G12").Value = 28.73869333333333
$ws.Range("H12").Value = 86.21608000000001
$ws.Range("I12").Value = 0.6207657201702246
$ws.Range("J12").Value = 0.6207657201702246
$ws.Range("M12").Value = 0.5637343333333333
$ws.Range("N12").Value = 1.691203
$ws.Range("O12").Value = 0.003475856580040611
$ws.Range("P12").Value = 0.003475856580040611
$ws.Range("Q12").Value = 16.20098812713778
$ws.Range("R12").Value = 145.80889314424
$ws.Range("S12").Value = 0.002157692613117324
$ws.Range("T12").Value = 0.002157692613117324

$ws.Range("G13").Value = 28.73869333333333
$ws.Range("H13").Value = 86.21608000000001
$ws.Range("I13").Value = 0.6207657201702246
$ws.Range("J13").Value = 0.6207657201702246
$ws.Range("M13").Value = 0.1137056666666667
$ws.Range("N13").Value = 0.341117
$ws.Range("O13").Value = 0.0007010830568617209
$ws.Range("P13").Value = 0.0007010830568617211
$ws.Range("Q13").Value = 3.267752284595555
$ws.Range("R13").Value = 29.40977056136
$ws.Range("S13").Value = 0.0004352083286919087
$ws.Range("T13").Value = 0.0004352083286919088

$ws.Range("G14").Value = 7.802604
$ws.Range("H14").Value = 23.407812
$ws.Range("I14").Value = 0.1685389462590879
$ws.Range("J14").Value = 0.1685389462590879
$ws.Range("M14").Value = 0.6692693333333334
$ws.Range("N14").Value = 2.007808
$ws.Range("O14").Value = 0.004126561180566838
$ws.Range("P14").Value = 0.004126561180566839
$ws.Range("Q14").Value = 5.222043577344
$ws.Range("R14").Value = 46.99839219609601
$ws.Range("S14").Value = 0.0006954862730463927
$ws.Range("T14").Value = 0.0006954862730463929

$ws.Range("G15").Value = 7.802604
$ws.Range("H15").Value = 23.407812
$ws.Range("I15").Value = 0.1685389462590879
$ws.Range("J15").Value = 0.1685389462590879
$ws.Range("O15").Value = 0.9916964991825307
$ws.Range("P15").Value = 0.9916964991825309
$ws.Range("Q15").Value = 1254.963178207212
$ws.Range("R15").Value = 11294.66860386491
$ws.Range("S15").Value = 0.1671394829810502
$ws.Range("T15").Value = 0.1671394829810502

$ws.Range("G16").Value = 7.802604
$ws.Range("H16").Value = 23.407812
$ws.Range("I16").Value = 0.1685389462590879
$ws.Range("J16").Value = 0.1685389462590879
$ws.Range("M16").Value = 0.5637343333333333
$ws.Range("N16").Value = 1.691203
$ws.Range("O16").Value = 0.003475856580040611
$ws.Range("P16").Value = 0.003475856580040611
$ws.Range("Q16").Value = 4.398595764204
$ws.Range("R16").Value = 39.587361877836
$ws.Range("S16").Value = 0.0005858172053477617
$ws.Range("T16").Value = 0.0005858172053477617

$ws.Range("G17").Value = 7.802604
$ws.Range("H17").Value = 23.407812
$ws.Range("I17").Value = 0.1685389462590879
$ws.Range("J17").Value = 0.1685389462590879
$ws.Range("M17").Value = 0.1137056666666667
$ws.Range("N17").Value = 0.341117
$ws.Range("O17").Value = 0.0007010830568617209
$ws.Range("P17").Value = 0.0007010830568617211
$ws.Range("Q17").Value = 0.8872002895559999
$ws.Range("R17").Value = 7.984802606004
$ws.Range("S17").Value = 0.0001181597996435747
$ws.Range("T17").Value = 0.0001181597996435747
